# PPP.xlsx update
#  - correctif probleme insertion matiere dans bd
#  - mise a jour semestre etudiant lors de l'inscription
#
# Effects:
#   1) The "Groupe" column uses shared text values "2-A"/"2-B"/"2-C" (semester 2)
#      which must become "1-A"/"1-B"/"1-C" (semester 1) wherever they occur.
#   2) Each student record's "Numero" (date-like id in column A) moves from the
#      2015 promotion to the 2017 promotion (+20000), and the computed
#      "Moyenne de l'etudiant" (column E) is refreshed to its new value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Update group/semester labels wherever they appear on the sheet ---
$groupMap = @{
    "2-B" = "1-B"
    "2-C" = "1-C"
    "2-A" = "1-A"
}

$used = $ws.UsedRange
foreach ($row in 1..$used.Rows.Count) {
    $cell = $ws.Cells.Item($row, 4)  # column D = "Groupe"
    $val = $cell.Value()
    if ($val -ne $null -and $groupMap.ContainsKey([string]$val)) {
        $cell.Value = $groupMap[[string]$val]
    }
}

# --- 2) Update per-row Numero (A) and Moyenne (E) ---
$rows = @(
    @{Row=3; A=20170926; E=19},
    @{Row=4; A=20170927; E=5},
    @{Row=5; A=20170928; E=9},
    @{Row=6; A=20170929; E=10},
    @{Row=7; A=20170930; E=18},
    @{Row=8; A=20170931; E=5},
    @{Row=9; A=20170932; E=14},
    @{Row=10; A=20170933; E=7},
    @{Row=11; A=20170934; E=19},
    @{Row=12; A=20170935; E=5},
    @{Row=13; A=20170936; E=9},
    @{Row=14; A=20170937; E=10},
    @{Row=15; A=20170938; E=10},
    @{Row=16; A=20170939; E=13},
    @{Row=17; A=20170940; E=15},
    @{Row=18; A=20170941; E=6},
    @{Row=19; A=20170942; E=18},
    @{Row=20; A=20170943; E=5},
    @{Row=21; A=20170944; E=19},
    @{Row=22; A=20170945; E=17},
    @{Row=23; A=20170946; E=14},
    @{Row=24; A=20170947},
    @{Row=25; A=20170948; E=6},
    @{Row=26; A=20170949; E=7},
    @{Row=27; A=20170950; E=18},
    @{Row=28; A=20170951},
    @{Row=29; A=20170952; E=14},
    @{Row=30; A=20170953; E=16},
    @{Row=31; A=20170954; E=10},
    @{Row=32; A=20170955; E=17},
    @{Row=33; A=20170956; E=5},
    @{Row=34; A=20170957; E=17},
    @{Row=35; A=20170958},
    @{Row=36; A=20170959; E=14},
    @{Row=37; A=20170960; E=11},
    @{Row=38; A=20170961; E=7},
    @{Row=39; A=20170962; E=19},
    @{Row=40; A=20170963; E=15},
    @{Row=41; A=20170964; E=13},
    @{Row=42; A=20170965; E=7},
    @{Row=43; A=20170966; E=13},
    @{Row=44; A=20170967; E=9},
    @{Row=45; A=20170968; E=14},
    @{Row=46; A=20170969; E=7},
    @{Row=47; A=20170970; E=18},
    @{Row=48; A=20170971; E=6},
    @{Row=49; A=20170972; E=13},
    @{Row=50; A=20170973; E=11},
    @{Row=51; A=20170974; E=9},
    @{Row=52; A=20170975; E=5},
    @{Row=53; A=20170976; E=8},
    @{Row=54; A=20170977; E=13},
    @{Row=55; A=20170978; E=13},
    @{Row=56; A=20170979; E=13},
    @{Row=57; A=20170980; E=13},
    @{Row=58; A=20170981; E=19},
    @{Row=59; A=20170982; E=7},
    @{Row=60; A=20170983; E=12},
    @{Row=61; A=20170984},
    @{Row=62; A=20170985; E=12},
    @{Row=63; A=20170986; E=12}
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    if ($r.ContainsKey("E")) {
        $ws.Cells.Item($r.Row, 5).Value = $r.E
    }
}

$wb.Save()
